$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date number format) from the cell above onto the new
# row's date cell, then add the new row of race-result data (row 21).
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A21").Value = 46073
$ws.Range("B21").Value = "Volta ao Algarve em Bicicleta"
$ws.Range("C21").Value = "Stage 3"
$ws.Range("D21").Value = "Filippo Ganna"
$ws.Range("E21").Value = "Juan Ayuso"
$ws.Range("F21").Value = "Jakob Söderqvist"
$ws.Range("G21").Value = "Paul Seixas"
$ws.Range("H21").Value = "Thymen Arensman"
$ws.Range("I21").Value = "Kévin Vauquelin"
$ws.Range("J21").Value = "Stefan Küng"
$ws.Range("K21").Value = "Héctor Álvarez"
$ws.Range("L21").Value = "Florian Lipowitz"
$ws.Range("M21").Value = "João Almeida"

$ws.Range("B32").Select()
